$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Repeat-Char([string]$ch, [int]$count) {
    $result = ""
    for ($i = 0; $i -lt $count; $i++) {
        $result = $result + $ch
    }
    return $result
}

# Build the "Jennifer Lighter..." author list string with a given number of
# spaces after each comma separator (mirrors the Springer API bug where the
# abstract/author refetch kept padding the separator with extra whitespace).
function Make-JenniferAuthors([int]$spaces) {
    $sep = "," + (Repeat-Char " " $spaces)
    $parts = @(
        "Jennifer%Lighter%Jennifer.Lighter@nyumc.org%1",
        "Michael%Phillips%NULL%1",
        "Sarah%Hochman%NULL%1",
        "Stephanie%Sterling%NULL%1",
        "Diane%Johnson%NULL%1",
        "Fritz%Francois%NULL%0",
        "Anna%Stachel%NULL%1"
    )
    return "[" + ($parts -join $sep) + "]"
}

# Build the "J. Wu..." author list string with a given number of spaces
# after each comma separator.
function Make-WuAuthors([int]$spaces) {
    $sep = "," + (Repeat-Char " " $spaces)
    $parts = @(
        "J.%Wu%xref no email%1",
        "W.%Li%xref no email%1",
        "X.%Shi%xref no email%1",
        "Z.%Chen%xref no email%1",
        "B.%Jiang%xref no email%1",
        "J.%Liu%xref no email%1",
        "D.%Wang%xref no email%1",
        "C.%Liu%xref no email%1",
        "Y.%Meng%xref no email%1",
        "L.%Cui%xref no email%1",
        "J.%Yu%xref no email%1",
        "H.%Cao%xref no email%1",
        "L.%Li%xref no email%1"
    )
    return "[" + ($parts -join $sep) + "]"
}

# The buggy refetch logic ran an extra round for each row, each time
# appending a new shared string entry before settling on the final value.
$ws.Range("E2").Value = Make-JenniferAuthors 9
$ws.Range("E3").Value = Make-WuAuthors 3

$ws.Range("E2").Value = Make-JenniferAuthors 10
$ws.Range("E3").Value = Make-WuAuthors 4
